# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (K instead of Strike#, recomputed std/mean, s_vals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 1
    14 = 3
    15 = 3
    16 = 3
    17 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
